$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "₹ 21,063"
$ws.Range("D2").Value = "₹ 134,802"
$ws.Range("C3").Value = "₹ 206,210"
$ws.Range("D3").Value = "₹ 1,301,801"
$ws.Range("C4").Value = "₹ 19,945"
$ws.Range("D4").Value = "₹ 96,167"
